$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final roster (column A), in top-to-bottom row order. "Buddy Hield" was
# dropped and replaced with "Miles McBride" (New York Knicks); the rest of
# the roster was re-entered in a new order.
$names = @(
    "Miles McBride",
    "Anfernee Simons",
    "Collin Sexton",
    "Fred VanVleet",
    "Paul George",
    "Giannis Antetokounmpo",
    "Anthony Edwards",
    "Aaron Gordon",
    "Jaren Jackson Jr.",
    "Ivica Zubac",
    "Khris Middleton",
    "Tyus Jones",
    "James Harden",
    "Jaden Ivey",
    "Jayson Tatum",
    "Zion Williamson",
    "Andrew Wiggins"
)

$positions = @(
    "PG,SG",
    "PG,SG",
    "PG,SG",
    "PG",
    "SG,SF,PF",
    "PF,C",
    "SG,SF",
    "PF,C",
    "PF,C",
    "C",
    "SF",
    "PG",
    "PG,SG",
    "PG,SG",
    "SF,PF",
    "PF,C",
    "SF,PF"
)

$teams = @(
    "New York Knicks",
    "Portland Trail Blazers",
    "Utah Jazz",
    "Houston Rockets",
    "Philadelphia 76ers",
    "Milwaukee Bucks",
    "Minnesota Timberwolves",
    "Denver Nuggets",
    "Memphis Grizzlies",
    "LA Clippers",
    "Milwaukee Bucks",
    "Phoenix Suns",
    "LA Clippers",
    "Detroit Pistons",
    "Boston Celtics",
    "New Orleans Pelicans",
    "Golden State Warriors"
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $names[$i]
}
for ($i = 0; $i -lt $positions.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $positions[$i]
}
for ($i = 0; $i -lt $teams.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $teams[$i]
}
